$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 301, shifting existing rows 301..400 down to 302..401.
$ws.Rows.Item(301).Insert()

# Populate the new row 301 with a duplicate of the (now shifted) row's data,
# but with an updated date (new weekly observation for the same price point).
$ws.Cells.Item(301, 1).Value  = 8
$ws.Cells.Item(301, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(301, 3).Value  = "Coquimbo"
$ws.Cells.Item(301, 4).Value  = 44524
$ws.Cells.Item(301, 5).Value  = 4
$ws.Cells.Item(301, 6).Value  = 100112043
$ws.Cells.Item(301, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(301, 8).Value  = "Sin especificar"
$ws.Cells.Item(301, 9).Value  = "Primera"
$ws.Cells.Item(301, 10).Value = 800
$ws.Cells.Item(301, 11).Value = 6500
$ws.Cells.Item(301, 12).Value = 7000
$ws.Cells.Item(301, 13).Value = 6750
$ws.Cells.Item(301, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(301, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(301, 16).Value = 112
$ws.Cells.Item(301, 17).Value = 60
$ws.Cells.Item(301, 18).Value = "Hortaliza"
